# MitsosBarton2006Ex312 - regenerate the non-convex experiment point
# (alpha = 0 stationary point), updating x, y and all dependent
# expression/evaluation/vector cells across the workbook's sheets.
#
# All of the affected cells store their (numeric-looking) content as
# literal TEXT, not as numbers (t="s" shared-string cells in the
# original file). A plain `Range.Value = "-1.8"` assignment would be
# auto-coerced to a real number by Excel, so instead each value is
# entered as a text-formula ("=""-1.8""") and then flattened back to a
# static value with Copy + PasteSpecial(xlPasteValues) - this keeps the
# cell's stored type as text without touching any cell's number format
# or style.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$text)
    $escaped = $text.Replace('"', '""')
    $range = $ws.Range($addr)
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

# NOTE: workbook has two sheets whose names differ only by case
# ("Vector_bf" / "Vector_BF") - Worksheets.Item(<name>) resolves
# case-insensitively and would collide, so sheets are addressed by
# their (1-based) tab position instead, matching xl/workbook.xml:
#   1 Funciones_Objetivo         5 Vector_bf
#   2 Restricciones_del_lider    6 Vector_BF
#   3 Restricciones_del_follower 7 Vector_Alpha
#   4 Punto_modificado

# --- Restricciones_del_lider --------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2 "A2" "0.8 - x"
Set-TextValue $ws2 "B2" "-1.8"
Set-TextValue $ws2 "D2" "0.74"
Set-TextValue $ws2 "A3" "-0.8 + x"
Set-TextValue $ws2 "B3" "-0.19999999999999996"
Set-TextValue $ws2 "D3" "0.96"

# --- Restricciones_del_follower -----------------------------------------
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3 "A2" "1.85 - y"
Set-TextValue $ws3 "B2" "-2.85"
Set-TextValue $ws3 "D2" "0.76"
Set-TextValue $ws3 "E2" "3.8"
Set-TextValue $ws3 "F2" "5.1"
Set-TextValue $ws3 "A3" "-1.85 + y"
Set-TextValue $ws3 "B3" "0.8500000000000001"
Set-TextValue $ws3 "D3" "0.2"
Set-TextValue $ws3 "E3" "0"
Set-TextValue $ws3 "F3" "3.7"

# --- Punto_modificado ------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4 "A2" "0.8"
Set-TextValue $ws4 "B2" "1.85"

# --- Vector_bf --------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
Set-TextValue $ws5 "A2" "-9.14325"

# --- Vector_BF --------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)
Set-TextValue $ws6 "A2" "-1.07"
Set-TextValue $ws6 "A3" "-34.0"
